# Aggiornamento gestione test K
# - Row 15 (ID 44, VALIDAZIONE_LAB_TIMEOUT): update expected-result note (col P)
#   to describe the timeout (504) behaviour instead of the stale "Referto
#   prodotto correttamente..." text.
# - Row 17 (ID 53, VALIDAZIONE_CDA2_LAB_CT7_KO): this row had been mistakenly
#   filled in with "OK case" style data (timestamp/traceid/workflow id and
#   SI/NO applicability flags). Re-align it with the other KO rows: clear the
#   stray values and record the real KO note in column K, with the
#   "GESTIONE ERRORE" flag (J) switched to NO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (touched first so its new shared string is interned before --
#     row 15's, matching the author's original shared-string ordering) ---
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = "NO"
$ws.Range("K17").Value = "L'applicativo è stato aggiornato in modo da salvare in automatico il CF in maiuscolo"
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("O17").Value = ""
$ws.Range("P17").Value = ""

# --- Row 15 -------------------------------------------------------------
$ws.Range("P15").Value = "In caso di timeout (errore 504) l'esecuzione prosegue ed il referto viene prodotto correttamente. Verrà visualizzato sulla richiesta l'esito negativo dell'invio a FSE e l'utente può effettuare il reinvio al momento oppure in un secondo momento"
